$wb = $excel.ActiveWorkbook

function Expand-ProductRows($ws) {
    # Insert a spacer row before each of the original data rows 9,10,11,12,13
    # (sequential insert points shift down by 1 each time: 9, 11, 13, 15, 17)
    $insertPoints = @(9, 11, 13, 15, 17)
    foreach ($p in $insertPoints) {
        $ws.Rows.Item($p).Insert()
        $srcRow = $p + 1
        $ws.Range("C" + $srcRow + ":E" + $srcRow).Copy()
        $ws.Range("C" + $p + ":E" + $p).PasteSpecial(-4122)
    }

    # Fill column B with "-" placeholder text for every row 9..18 except row 14
    # (row 14 keeps the existing "Product Group" header label already present there)
    for ($r = 9; $r -le 18; $r++) {
        if ($r -ne 14) {
            $ws.Range("B" + $r).Value = "-"
        }
    }
}

$ws1 = $wb.Worksheets.Item("Plan1")
$ws2 = $wb.Worksheets.Item("Plan2")

Expand-ProductRows $ws1
Expand-ProductRows $ws2

# Update named ranges that referenced the old C8:C13 / D8:D13 / E8:E13 spans
$wb.Names.Item(5).RefersTo = "=Plan2!`$C`$8:`$C`$18"
$wb.Names.Item(6).RefersTo = "=Plan1!`$C`$8:`$C`$18"
$wb.Names.Item(7).RefersTo = "=Plan2!`$D`$8:`$D`$18"
$wb.Names.Item(8).RefersTo = "=Plan1!`$D`$8:`$D`$18"
$wb.Names.Item(9).RefersTo = "=Plan2!`$E`$8:`$E`$18"
$wb.Names.Item(10).RefersTo = "=Plan1!`$E`$8:`$E`$18"

# Update the selections to reflect where the user ended up after the edit
$ws2.Activate() | Out-Null
$ws2.Range("B19").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B19").Select() | Out-Null
